$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 327, pushing existing rows 327:396 down to 328:397.
$ws.Rows.Item(327).Insert()

# Populate the newly inserted row 327 with the new record.
$ws.Cells.Item(327, 1).Value = 4
$ws.Cells.Item(327, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(327, 3).Value = "Los Lagos"
$ws.Cells.Item(327, 4).Value = 44641
$ws.Cells.Item(327, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(327, 5).Value = 10
$ws.Cells.Item(327, 6).Value = 100112006
$ws.Cells.Item(327, 7).Value = "Repollo"
$ws.Cells.Item(327, 8).Value = "Copenhague"
$ws.Cells.Item(327, 9).Value = "Primera"
$ws.Cells.Item(327, 10).Value = 500
$ws.Cells.Item(327, 11).Value = 2000
$ws.Cells.Item(327, 12).Value = 2000
$ws.Cells.Item(327, 13).Value = 2000
$ws.Cells.Item(327, 14).Value = "$/unidad"
$ws.Cells.Item(327, 15).Value = "Región del Maule"
$ws.Cells.Item(327, 16).Value = 2000
$ws.Cells.Item(327, 17).Value = 1
$ws.Cells.Item(327, 18).Value = "Hortaliza"
